# Auto-generated edit script: updates Leve market-price derived columns
# (currentAveragePrice / LevePrice / LeveProfit) across multiple job sheets
# to reflect refreshed market data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Row 15 (ALC): Morning Glass of Ether | Ether
$ws_ALC.Range("H15").Value = 1197.0533
$ws_ALC.Range("I15").Value = 1197.0533
$ws_ALC.Range("K15").Value = 3591.1599
$ws_ALC.Range("M15").Value = -3422.1599

# Row 17 (ALC): One for the Road | Potion
$ws_ALC.Range("H17").Value = 403062.44
$ws_ALC.Range("J17").Value = 403062.44
$ws_ALC.Range("L17").Value = 1209187.32
$ws_ALC.Range("N17").Value = -1209523.32

# Row 43 (ALC): Growing Is Knowing | Growth Formula Gamma
$ws_ALC.Range("H43").Value = 4060667.5
$ws_ALC.Range("J43").Value = 18443.428
$ws_ALC.Range("L43").Value = 18443.428
$ws_ALC.Range("N43").Value = -18581.428

# Row 70 (ALC): Consecrating Congregation | Holy Water
$ws_ALC.Range("H70").Value = 71345
$ws_ALC.Range("I70").Value = 144314.28
$ws_ALC.Range("K70").Value = 432942.84
$ws_ALC.Range("M70").Value = -432672.84

# Row 73 (ALC): Curbing the Contagion (L) | Holy Water
$ws_ALC.Range("H73").Value = 71345
$ws_ALC.Range("I73").Value = 144314.28
$ws_ALC.Range("K73").Value = 432942.84
$ws_ALC.Range("M73").Value = -432006.84

# Row 112 (ALC): Making Ends Meet | Superior Spiritbond Potion
$ws_ALC.Range("H112").Value = 627809.0600000001
$ws_ALC.Range("I112").Value = 1461.3334
$ws_ALC.Range("J112").Value = 772350.9
$ws_ALC.Range("K112").Value = 4384.0002
$ws_ALC.Range("L112").Value = 2317052.7
$ws_ALC.Range("M112").Value = -3276.0002
$ws_ALC.Range("N112").Value = -2319268.7

# Row 135 (ALC): For Tired Minds | Grade 1 Gemsap of Intelligence
$ws_ALC.Range("H135").Value = 3026.186
$ws_ALC.Range("I135").Value = 2068.2424
$ws_ALC.Range("J135").Value = 6187.4
$ws_ALC.Range("K135").Value = 18614.1816
$ws_ALC.Range("L135").Value = 55686.6
$ws_ALC.Range("M135").Value = -16079.1816
$ws_ALC.Range("N135").Value = -60756.6

# Row 137 (ALC): Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws_ALC.Range("H137").Value = 4110.615
$ws_ALC.Range("I137").Value = 3451.4443
$ws_ALC.Range("J137").Value = 5593.75
$ws_ALC.Range("K137").Value = 10354.3329
$ws_ALC.Range("L137").Value = 16781.25
$ws_ALC.Range("M137").Value = -7804.332900000001
$ws_ALC.Range("N137").Value = -21881.25

# Row 32 (ARM): Ingot We Trust | Steel Ingot
$ws_ARM.Range("H32").Value = 2820.2917
$ws_ARM.Range("I32").Value = 2749.9785
$ws_ARM.Range("K32").Value = 2749.9785
$ws_ARM.Range("M32").Value = -2462.9785

# Row 88 (ARM): The Mast Chance | Adamantite Rivets
$ws_ARM.Range("H88").Value = 1230
$ws_ARM.Range("I88").Value = 547.75
$ws_ARM.Range("J88").Value = 2139.6667
$ws_ARM.Range("K88").Value = 547.75
$ws_ARM.Range("L88").Value = 2139.6667
$ws_ARM.Range("M88").Value = -141.75
$ws_ARM.Range("N88").Value = -2951.6667

# Row 91 (ARM): The Rose and the Riveter (L) | Adamantite Rivets
$ws_ARM.Range("H91").Value = 1230
$ws_ARM.Range("I91").Value = 547.75
$ws_ARM.Range("J91").Value = 2139.6667
$ws_ARM.Range("K91").Value = 547.75
$ws_ARM.Range("L91").Value = 2139.6667
$ws_ARM.Range("M91").Value = 856.25
$ws_ARM.Range("N91").Value = -4947.6667

# Row 97 (ARM): Ore for Me | High Steel Ingot
$ws_ARM.Range("H97").Value = 1602.2354
$ws_ARM.Range("I97").Value = 898.13794
$ws_ARM.Range("K97").Value = 898.13794
$ws_ARM.Range("M97").Value = -402.13794

# Row 112 (ARM): Wrapped Knuckles | Deepgold Gloves of Fending
$ws_ARM.Range("H112").Value = 0
$ws_ARM.Range("J112").Value = 0
$ws_ARM.Range("N112").ClearContents()

# Row 133 (ARM): Shielding My Students | Mountain Chromite Tower Shield
$ws_ARM.Range("H133").Value = 106969
$ws_ARM.Range("J133").Value = 106969
$ws_ARM.Range("L133").Value = 106969
$ws_ARM.Range("N133").Value = -112029

# Row 107 (BSM): The Gold Experience | Deepgold Nugget
$ws_BSM.Range("H107").Value = 1170.5238
$ws_BSM.Range("I107").Value = 1004.2632
$ws_BSM.Range("K107").Value = 1004.2632
$ws_BSM.Range("M107").Value = 915.7368

# Row 134 (BSM): Ruthenium Supremium | Ruthenium Ingot
$ws_BSM.Range("H134").Value = 3873.7593
$ws_BSM.Range("I134").Value = 1414.7297
$ws_BSM.Range("J134").Value = 9225.764999999999
$ws_BSM.Range("K134").Value = 4244.189100000001
$ws_BSM.Range("L134").Value = 27677.295
$ws_BSM.Range("M134").Value = -1709.189100000001
$ws_BSM.Range("N134").Value = -32747.295

# Row 31 (CRP): Wall Not Found | Walnut Lumber
$ws_CRP.Range("H31").Value = 2236.0962
$ws_CRP.Range("I31").Value = 1655.3096
$ws_CRP.Range("J31").Value = 4675.4
$ws_CRP.Range("K31").Value = 1655.3096
$ws_CRP.Range("L31").Value = 4675.4
$ws_CRP.Range("M31").Value = -1360.3096
$ws_CRP.Range("N31").Value = -5265.4

# Row 34 (CRP): Armoires of the Rich and Famous | Walnut Lumber
$ws_CRP.Range("H34").Value = 2236.0962
$ws_CRP.Range("I34").Value = 1655.3096
$ws_CRP.Range("J34").Value = 4675.4
$ws_CRP.Range("K34").Value = 1655.3096
$ws_CRP.Range("L34").Value = 4675.4
$ws_CRP.Range("M34").Value = -1453.3096
$ws_CRP.Range("N34").Value = -5079.4

# Row 58 (CRP): You Do the Heavy Lifting | Mahogany Lumber
$ws_CRP.Range("H58").Value = 1113.5
$ws_CRP.Range("I58").Value = 1113.5
$ws_CRP.Range("K58").Value = 1113.5
$ws_CRP.Range("M58").Value = -910.5

# Row 105 (CRP): Zelkova, My Love | Zelkova Lumber
$ws_CRP.Range("H105").Value = 861.25
$ws_CRP.Range("I105").Value = 898.3333
$ws_CRP.Range("K105").Value = 898.3333
$ws_CRP.Range("M105").Value = 848.6667

# Row 132 (CRP): Hull Lotta Damage | Ginseng Lumber
$ws_CRP.Range("H132").Value = 1440.8906
$ws_CRP.Range("I132").Value = 1336.95
$ws_CRP.Range("J132").Value = 3000
$ws_CRP.Range("K132").Value = 4010.85
$ws_CRP.Range("L132").Value = 9000
$ws_CRP.Range("M132").Value = -1480.85
$ws_CRP.Range("N132").Value = -14060

# Row 134 (CRP): Wood You Be Quiet | Ceiba Lumber
$ws_CRP.Range("H134").Value = 2163.6482
$ws_CRP.Range("I134").Value = 1842.1459
$ws_CRP.Range("K134").Value = 5526.4377
$ws_CRP.Range("M134").Value = -2991.4377

# Row 136 (CRP): Turali Quality | Dark Mahogany Lumber
$ws_CRP.Range("H136").Value = 1113.5
$ws_CRP.Range("I136").Value = 1113.5
$ws_CRP.Range("K136").Value = 3340.5
$ws_CRP.Range("M136").Value = -790.5

# Row 56 (CUL): Culture Club | Crowned Pie
$ws_CUL.Range("H56").Value = 959810.9399999999
$ws_CUL.Range("I56").Value = 959810.9399999999
$ws_CUL.Range("K56").Value = 959810.9399999999
$ws_CUL.Range("M56").Value = -959280.9399999999

# Row 131 (CUL): The Mountain Steeped | Tsai tou Vounou
$ws_CUL.Range("H131").Value = 1224.2222
$ws_CUL.Range("I131").Value = 919.8
$ws_CUL.Range("J131").Value = 1604.75
$ws_CUL.Range("K131").Value = 2759.4
$ws_CUL.Range("L131").Value = 4814.25
$ws_CUL.Range("M131").Value = 2280.6
$ws_CUL.Range("N131").Value = -14894.25

# Row 137 (CUL): Creative Chocolate | Gateau au Chocolat
$ws_CUL.Range("H137").Value = 6524.9165
$ws_CUL.Range("J137").Value = 7837.5
$ws_CUL.Range("L137").Value = 23512.5
$ws_CUL.Range("N137").Value = -33712.5

# Row 100 (GSM): Hair-raising Action | Durium Hairpin of Maiming
$ws_GSM.Range("H100").Value = 0
$ws_GSM.Range("J100").Value = 0
$ws_GSM.Range("N100").ClearContents()

# Row 113 (GSM): Copious Crystal Cannons | Manasilver Nugget
$ws_GSM.Range("H113").Value = 1653.1666
$ws_GSM.Range("I113").Value = 1398.5
$ws_GSM.Range("K113").Value = 1398.5
$ws_GSM.Range("M113").Value = 771.5

# Row 132 (GSM): On Board for Lar | Lar Ingot
$ws_GSM.Range("H132").Value = 7421.3887
$ws_GSM.Range("I132").Value = 7769.706
$ws_GSM.Range("J132").Value = 1500
$ws_GSM.Range("K132").Value = 23309.118
$ws_GSM.Range("L132").Value = 4500
$ws_GSM.Range("M132").Value = -20779.118
$ws_GSM.Range("N132").Value = -9560

# Row 42 (LTW): Slave to Fashion | Boarskin Choker
$ws_LTW.Range("H42").Value = 13407.143
$ws_LTW.Range("I42").Value = 8691.333000000001
$ws_LTW.Range("J42").Value = 16944
$ws_LTW.Range("K42").Value = 8691.333000000001
$ws_LTW.Range("L42").Value = 16944
$ws_LTW.Range("M42").Value = -8128.333000000001
$ws_LTW.Range("N42").Value = -18070

# Row 46 (LTW): Supply Side Logic | Boar Leather
$ws_LTW.Range("H46").Value = 11344.77
$ws_LTW.Range("I46").Value = 3179
$ws_LTW.Range("J46").Value = 24410
$ws_LTW.Range("K46").Value = 3179
$ws_LTW.Range("L46").Value = 24410
$ws_LTW.Range("M46").Value = -2991
$ws_LTW.Range("N46").Value = -24786

# Row 49 (LTW): First They Came for the Heretics | Boarskin Choker
$ws_LTW.Range("H49").Value = 13407.143
$ws_LTW.Range("I49").Value = 8691.333000000001
$ws_LTW.Range("J49").Value = 16944
$ws_LTW.Range("K49").Value = 8691.333000000001
$ws_LTW.Range("L49").Value = 16944
$ws_LTW.Range("M49").Value = -8544.333000000001
$ws_LTW.Range("N49").Value = -17238

# Row 61 (LTW): Spelling Me Softly | Raptor Leather
$ws_LTW.Range("H61").Value = 1662.2778
$ws_LTW.Range("I61").Value = 1557.625
$ws_LTW.Range("K61").Value = 1557.625
$ws_LTW.Range("M61").Value = -1355.625

# Row 68 (LTW): You Could Say It's a Moving Target | Wyvern Leather
$ws_LTW.Range("H68").Value = 7174.0835
$ws_LTW.Range("I68").Value = 6022.5
$ws_LTW.Range("J68").Value = 7749.875
$ws_LTW.Range("K68").Value = 6022.5
$ws_LTW.Range("L68").Value = 7749.875
$ws_LTW.Range("M68").Value = -5273.5
$ws_LTW.Range("N68").Value = -9247.875

# Row 71 (LTW): They Call It Bloody Mary (L) | Wyvern Leather
$ws_LTW.Range("H71").Value = 7174.0835
$ws_LTW.Range("I71").Value = 6022.5
$ws_LTW.Range("J71").Value = 7749.875
$ws_LTW.Range("K71").Value = 30112.5
$ws_LTW.Range("L71").Value = 38749.375
$ws_LTW.Range("M71").Value = -26368.5
$ws_LTW.Range("N71").Value = -46237.375

# Row 93 (LTW): Hide to Go Seek | Gagana Leather
$ws_LTW.Range("H93").Value = 4921.5
$ws_LTW.Range("I93").Value = 4925.8
$ws_LTW.Range("J93").Value = 4900
$ws_LTW.Range("K93").Value = 4925.8
$ws_LTW.Range("L93").Value = 4900
$ws_LTW.Range("M93").Value = -3677.8
$ws_LTW.Range("N93").Value = -7396

# Row 112 (LTW): A Slippery Slope | Gliderskin Boots of Casting
$ws_LTW.Range("H112").Value = 71313.336
$ws_LTW.Range("J112").Value = 71313.336
$ws_LTW.Range("L112").Value = 71313.336
$ws_LTW.Range("N112").Value = -74267.336

# Row 113 (LTW): Peace in Rest | Atrociraptor Leather
$ws_LTW.Range("H113").Value = 1662.2778
$ws_LTW.Range("I113").Value = 1557.625
$ws_LTW.Range("K113").Value = 1557.625
$ws_LTW.Range("M113").Value = 612.375

# Row 119 (LTW): Fit for a Friend | Swallowskin Gloves of Fending
$ws_LTW.Range("H119").Value = 0
$ws_LTW.Range("J119").Value = 0
$ws_LTW.Range("N119").ClearContents()

# Row 132 (LTW): Tenets of Tanning | Silver Lobo Leather
$ws_LTW.Range("H132").Value = 1668.0638
$ws_LTW.Range("I132").Value = 1241.4634
$ws_LTW.Range("J132").Value = 4583.1665
$ws_LTW.Range("K132").Value = 3724.3902
$ws_LTW.Range("L132").Value = 13749.4995
$ws_LTW.Range("M132").Value = -1194.3902
$ws_LTW.Range("N132").Value = -18809.4995

# Row 62 (WVR): Pride Up in Smoke | Rainbow Cloth
$ws_WVR.Range("H62").Value = 8401502
$ws_WVR.Range("I62").Value = 90025.44500000001
$ws_WVR.Range("J62").Value = 33335930
$ws_WVR.Range("K62").Value = 90025.44500000001
$ws_WVR.Range("L62").Value = 33335930
$ws_WVR.Range("M62").Value = -89401.44500000001
$ws_WVR.Range("N62").Value = -33337178

# Row 65 (WVR): Desperate for Diversionaries (L) | Rainbow Cloth
$ws_WVR.Range("H65").Value = 8401502
$ws_WVR.Range("I65").Value = 90025.44500000001
$ws_WVR.Range("J65").Value = 33335930
$ws_WVR.Range("K65").Value = 450127.225
$ws_WVR.Range("L65").Value = 166679650
$ws_WVR.Range("M65").Value = -447007.225
$ws_WVR.Range("N65").Value = -166685890

# Row 81 (WVR): Where the Dragonflies, the Net Catches | Crawler Silk
$ws_WVR.Range("H81").Value = 43433.11
$ws_WVR.Range("I81").Value = 94415
$ws_WVR.Range("J81").Value = 8383.0625
$ws_WVR.Range("K81").Value = 188830
$ws_WVR.Range("L81").Value = 16766.125
$ws_WVR.Range("M81").Value = -187769
$ws_WVR.Range("N81").Value = -18888.125

# Row 84 (WVR): To Kill a Dragon on Nameday (L) | Crawler Silk
$ws_WVR.Range("H84").Value = 43433.11
$ws_WVR.Range("I84").Value = 94415
$ws_WVR.Range("J84").Value = 8383.0625
$ws_WVR.Range("K84").Value = 944150
$ws_WVR.Range("L84").Value = 83830.625
$ws_WVR.Range("M84").Value = -938846
$ws_WVR.Range("N84").Value = -94438.625

# Row 92 (WVR): Modest Beginnings | Bloodhempen Culottes of Casting
$ws_WVR.Range("H92").Value = 0
$ws_WVR.Range("J92").Value = 0
$ws_WVR.Range("N92").ClearContents()

# Row 100 (WVR): Of Great Import | Kudzu Thread
$ws_WVR.Range("H100").Value = 991.6585
$ws_WVR.Range("I100").Value = 1143.2593
$ws_WVR.Range("J100").Value = 699.2857
$ws_WVR.Range("K100").Value = 2286.5186
$ws_WVR.Range("L100").Value = 1398.5714
$ws_WVR.Range("M100").Value = -1745.5186
$ws_WVR.Range("N100").Value = -2480.5714

# Row 132 (WVR): Comfy Cabins | Snow Cotton Cloth
$ws_WVR.Range("H132").Value = 1363.3387
$ws_WVR.Range("I132").Value = 1259.0862
$ws_WVR.Range("K132").Value = 3777.2586
$ws_WVR.Range("M132").Value = -1247.2586

# Row 136 (WVR): Weaving the Envelope | Sarcenet Cloth
$ws_WVR.Range("H136").Value = 8718.104499999999
$ws_WVR.Range("I136").Value = 9776.081
$ws_WVR.Range("J136").Value = 5159.4546
$ws_WVR.Range("K136").Value = 29328.243
$ws_WVR.Range("L136").Value = 15478.3638
$ws_WVR.Range("M136").Value = -26778.243
$ws_WVR.Range("N136").Value = -20578.3638
